$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the K2 formula (add missing quotes around C2/D2) and fill down to K63
$ws.Range("K2").Formula = '=CONCATENATE("INSERT INTO ARTICULOS VALUES ( ''",C2, "'', ''", D2, "'', ''", E2, "'', ''", F2, "'', ''", G2, "'', ''", H2, "'', ''", I2, "'', ''", J2, "'');")'
$ws.Range("K2:K63").FillDown()

# Column K width
$ws.Range("K1").ColumnWidth = 15.28515625

# Selection / view state
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("L12").Select()
